$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.200.95"
$ws.Range("E2").Value = "  +1.88%  "
$ws.Range("D3").Value = "3.386.37"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'586.86"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("D6").Value = "'180.15"
$ws.Range("E6").Value = "  +1.65%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.595"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").Value = "'0.195"
$ws.Range("E9").Value = "  +6.99%  "
$ws.Range("D10").Value = "'0.590"
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("D11").Value = "'48.57"
$ws.Range("E11").Value = "  +3.72%  "
$ws.Range("D12").Value = "'0.0000281"
$ws.Range("E12").Value = "  +2.76%  "
$ws.Range("D13").Value = "'676.60"
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").Value = "'8.62"
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("D15").Value = "3.933.95"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "69.243.03"
$ws.Range("E16").Value = "  +1.85%  "
$ws.Range("D17").Value = "3.389.90"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("E18").Value = "  +1.88%  "
$ws.Range("D19").Value = "'17.67"
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("D20").Value = "'11.34"
$ws.Range("E20").Value = "  +2.43%  "
$ws.Range("D21").Value = "'0.901"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").Value = "'17.07"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "'103.38"
$ws.Range("E24").Value = "  +4.54%  "
$ws.Range("D25").Value = "'3.91"
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("D27").Value = "'9.58"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("D28").Value = "'34.10"
$ws.Range("E28").Value = "  +3.69%  "
$ws.Range("D29").Value = "'8.70"
$ws.Range("E29").Value = "  +1.53%  "
$ws.Range("D30").Value = "'6.98"
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("D31").Value = "'3.71"
$ws.Range("E31").Value = "  +13.06%  "
$ws.Range("D32").Value = "'11.18"
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("D33").Value = "'553.87"
$ws.Range("E33").Value = "  -3.70%  "
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("D35").Value = "'57.97"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "3.692.10"
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("E38").Value = "  +6.20%  "
$ws.Range("D39").Value = "'35.11"
$ws.Range("E39").Value = "  +1.33%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'3.24"
$ws.Range("E40").Value = "  +1.38%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0704"
$ws.Range("E41").Value = "  +4.11%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "'0.338"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("D44").Value = "'0.0422"
$ws.Range("E44").Value = "  +4.01%  "
$ws.Range("D45").Value = "'3.27"
$ws.Range("E45").Value = "  -3.35%  "
$ws.Range("D46").Value = "'2.65"
$ws.Range("E46").Value = "  -0.55%  "
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("E48").Value = "  +5.77%  "
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "'131.77"
$ws.Range("E50").Value = "  +1.78%  "
$ws.Range("D51").Value = "'2.60"
$ws.Range("E51").Value = "  -1.32%  "
